$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Translate Lithuanian descriptions to English (and fix B2 to match the
# variable-name column, as in the source edit).
$ws.Range("B2").Value  = "hh_ident"
$ws.Range("B18").Value = "gender"
$ws.Range("B20").Value = "age"
$ws.Range("B21").Value = "marital"
$ws.Range("B32").Value = "status_in_house"
$ws.Range("B38").Value = "education"
$ws.Range("B49").Value = "employment"
$ws.Range("B67").Value = "employment_type"
$ws.Range("B72").Value = "job_contract"

# Restore the active-cell selection as captured in the saved file.
$ws.Range("F39").Select()
